$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Service ID values refreshed to the latest test run (8th Sept 2023)
$ws.Range("C2").Value  = "10118041"
$ws.Range("C3").Value  = "10118046"
$ws.Range("C4").Value  = "10117974"
$ws.Range("C5").Value  = "10118050"
$ws.Range("C6").Value  = "10118059"
$ws.Range("C7").Value  = "10118070"
$ws.Range("C8").Value  = "10118077"
$ws.Range("C9").Value  = "10118081"
$ws.Range("C10").Value = "10118088"
$ws.Range("C11").Value = "10118108"
$ws.Range("C12").Value = "10118129"
$ws.Range("C13").Value = "10118152"
$ws.Range("C14").Value = "10118173"
$ws.Range("C20").Value = "794659397220"
$ws.Range("C24").Value = "135232141"

# Result columns updated to PASS for the re-run rows
$ws.Range("E14").Value = "PASS"
$ws.Range("E17").Value = "PASS"
$ws.Range("E23").Value = "PASS"

# Fail-log text captured for one of the still-failing checks (kept in the
# shared string table even though no cell currently points at it, matching
# the source workbook's sharedStrings.xml).
$errorText = "no such element: Unable to locate element: {`"method`":`"css selector`",`"selector`":`"#lblServiceID`"}" + "`r`n" + `
"  (Session info: headless chrome=116.0.5845.142)" + "`r`n" + `
"For documentation on this error, please visit: https://selenium.dev/exceptions/#no_such_element" + "`r`n" + `
"Build info: version: '4.9.0', revision: 'd7057100a6'" + "`r`n" + `
"System info: os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '19.0.1'" + "`r`n" + `
"Driver info: org.openqa.selenium.chrome.ChromeDriver" + "`r`n" + `
"Command: [4fb340dd9329317d63e252b8a77319a0, findElement {using=id, value=lblServiceID}]" + "`r`n" + `
"Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 116.0.5845.142, chrome: {chromedriverVersion: 116.0.5845.96 (1a3918166880..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:52521}, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: WINDOWS, proxy: Proxy(), se:cdp: ws://localhost:52521/devtoo..., se:cdpVersion: 116.0.5845.142, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}" + "`r`n" + `
"Session ID: 4fb340dd9329317d63e252b8a77319a0"

$ws.Range("Z1").Value = $errorText
$ws.Range("Z1").Value = ""
